$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 40750
$ws.Range("J93").Value = 40750
$ws.Range("L93").Value = 40750
$ws.Range("N93").Value = -45742

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1266.8667
$ws.Range("I98").Value = 1009.56525
$ws.Range("J98").Value = 2112.2856
$ws.Range("K98").Value = 1009.56525
$ws.Range("L98").Value = 2112.2856
$ws.Range("M98").Value = 488.43475
$ws.Range("N98").Value = -5108.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3106647.2
$ws.Range("J112").Value = 3572499.2
$ws.Range("L112").Value = 10717497.6
$ws.Range("N112").Value = -10719713.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1266.8667
$ws.Range("I122").Value = 1009.56525
$ws.Range("J122").Value = 2112.2856
$ws.Range("K122").Value = 3028.69575
$ws.Range("L122").Value = 6336.8568
$ws.Range("M122").Value = -578.6957499999999
$ws.Range("N122").Value = -11236.8568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1056.8148
$ws.Range("J129").Value = 1161.439
$ws.Range("L129").Value = 3484.317
$ws.Range("N129").Value = -13484.317

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1383.9736
$ws.Range("I132").Value = 1046.2188
$ws.Range("J132").Value = 3185.3333
$ws.Range("K132").Value = 3138.6564
$ws.Range("L132").Value = 9555.999899999999
$ws.Range("M132").Value = -608.6564000000003
$ws.Range("N132").Value = -14615.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1826.25
$ws.Range("I137").Value = 1391.2222
$ws.Range("J137").Value = 2609.3
$ws.Range("K137").Value = 4173.6666
$ws.Range("L137").Value = 7827.900000000001
$ws.Range("M137").Value = -1623.6666
$ws.Range("N137").Value = -12927.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 634.25
$ws.Range("I2").Value = 666.24
$ws.Range("J2").Value = 367.66666
$ws.Range("K2").Value = 666.24
$ws.Range("L2").Value = 367.66666
$ws.Range("M2").Value = -553.24
$ws.Range("N2").Value = -593.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5969.46
$ws.Range("I32").Value = 4571.291
$ws.Range("K32").Value = 4571.291
$ws.Range("M32").Value = -4284.291

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 831.6667
$ws.Range("I110").Value = 698
$ws.Range("K110").Value = 698
$ws.Range("M110").Value = 1347

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 634.25
$ws.Range("I116").Value = 666.24
$ws.Range("J116").Value = 367.66666
$ws.Range("K116").Value = 666.24
$ws.Range("L116").Value = 367.66666
$ws.Range("M116").Value = 1627.76
$ws.Range("N116").Value = -4955.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4274689.5
$ws.Range("I122").Value = 5129227
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 15387681
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -15385231
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2780.239
$ws.Range("I132").Value = 1797.8286
$ws.Range("K132").Value = 5393.4858
$ws.Range("M132").Value = -2863.4858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 35787
$ws.Range("J133").Value = 35787
$ws.Range("L133").Value = 35787
$ws.Range("N133").Value = -40847

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 42959.5
$ws.Range("J139").Value = 42959.5
$ws.Range("L139").Value = 42959.5
$ws.Range("N139").Value = -53239.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 634.25
$ws.Range("I3").Value = 666.24
$ws.Range("J3").Value = 367.66666
$ws.Range("K3").Value = 666.24
$ws.Range("L3").Value = 367.66666
$ws.Range("M3").Value = -552.24
$ws.Range("N3").Value = -595.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 76924310
$ws.Range("I107").Value = 166667740
$ws.Range("K107").Value = 166667740
$ws.Range("M107").Value = -166665820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5708.387
$ws.Range("I134").Value = 8828.866
$ws.Range("J134").Value = 2782.9375
$ws.Range("K134").Value = 26486.598
$ws.Range("L134").Value = 8348.8125
$ws.Range("M134").Value = -23951.598
$ws.Range("N134").Value = -13418.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1322.9584
$ws.Range("I16").Value = 1258.1818
$ws.Range("J16").Value = 1377.7693
$ws.Range("K16").Value = 1258.1818
$ws.Range("L16").Value = 1377.7693
$ws.Range("M16").Value = -971.1818000000001
$ws.Range("N16").Value = -1951.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 355182.78
$ws.Range("I31").Value = 1880.4517
$ws.Range("J31").Value = 1197672.9
$ws.Range("K31").Value = 1880.4517
$ws.Range("L31").Value = 1197672.9
$ws.Range("M31").Value = -1585.4517
$ws.Range("N31").Value = -1198262.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 355182.78
$ws.Range("I34").Value = 1880.4517
$ws.Range("J34").Value = 1197672.9
$ws.Range("K34").Value = 1880.4517
$ws.Range("L34").Value = 1197672.9
$ws.Range("M34").Value = -1678.4517
$ws.Range("N34").Value = -1198076.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1801.0344
$ws.Range("I105").Value = 1922.381
$ws.Range("K105").Value = 1922.381
$ws.Range("M105").Value = -175.3810000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1322.9584
$ws.Range("I113").Value = 1258.1818
$ws.Range("J113").Value = 1377.7693
$ws.Range("K113").Value = 1258.1818
$ws.Range("L113").Value = 1377.7693
$ws.Range("M113").Value = 911.8181999999999
$ws.Range("N113").Value = -5717.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1887.5424
$ws.Range("I134").Value = 2245.4102
$ws.Range("K134").Value = 6736.230599999999
$ws.Range("M134").Value = -4201.230599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 9376.973
$ws.Range("I134").Value = 9719.076999999999
$ws.Range("J134").Value = 9191.666999999999
$ws.Range("K134").Value = 29157.231
$ws.Range("L134").Value = 27575.001
$ws.Range("M134").Value = -24087.231
$ws.Range("N134").Value = -37715.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 50001210
$ws.Range("I113").Value = 125000984
$ws.Range("J113").Value = 1353.3334
$ws.Range("K113").Value = 125000984
$ws.Range("L113").Value = 1353.3334
$ws.Range("M113").Value = -124998814
$ws.Range("N113").Value = -5693.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 36719624
$ws.Range("I122").Value = 66552456
$ws.Range("J122").Value = 2290.8462
$ws.Range("K122").Value = 199657368
$ws.Range("L122").Value = 6872.5386
$ws.Range("M122").Value = -199654918
$ws.Range("N122").Value = -11772.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 30997.234
$ws.Range("I7").Value = 45016.78
$ws.Range("J7").Value = 1683.6364
$ws.Range("K7").Value = 45016.78
$ws.Range("L7").Value = 1683.6364
$ws.Range("M7").Value = -44904.78
$ws.Range("N7").Value = -1907.6364

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 45457260
$ws.Range("I40").Value = 55557212
$ws.Range("K40").Value = 55557212
$ws.Range("M40").Value = -55557076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1317.3572
$ws.Range("I61").Value = 1410.3529
$ws.Range("J61").Value = 1173.6364
$ws.Range("K61").Value = 1410.3529
$ws.Range("L61").Value = 1173.6364
$ws.Range("M61").Value = -1208.3529
$ws.Range("N61").Value = -1577.6364

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1317.3572
$ws.Range("I113").Value = 1410.3529
$ws.Range("J113").Value = 1173.6364
$ws.Range("K113").Value = 1410.3529
$ws.Range("L113").Value = 1173.6364
$ws.Range("M113").Value = 759.6470999999999
$ws.Range("N113").Value = -5513.6364

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2629467
$ws.Range("I122").Value = 3108412
$ws.Range("J122").Value = 1252500
$ws.Range("K122").Value = 9325236
$ws.Range("L122").Value = 3757500
$ws.Range("M122").Value = -9322786
$ws.Range("N122").Value = -3762400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 30997.234
$ws.Range("I126").Value = 45016.78
$ws.Range("J126").Value = 1683.6364
$ws.Range("K126").Value = 135050.34
$ws.Range("L126").Value = 5050.9092
$ws.Range("M126").Value = -132580.34
$ws.Range("N126").Value = -9990.9092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10908872
$ws.Range("I132").Value = 13362423
$ws.Range("K132").Value = 40087269
$ws.Range("M132").Value = -40084739

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7936.7334
$ws.Range("I136").Value = 5023.4443
$ws.Range("J136").Value = 19589.889
$ws.Range("K136").Value = 15070.3329
$ws.Range("L136").Value = 58769.667
$ws.Range("M136").Value = -12520.3329
$ws.Range("N136").Value = -63869.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2206.7715
$ws.Range("I122").Value = 2155.2693
$ws.Range("J122").Value = 2355.5557
$ws.Range("K122").Value = 6465.8079
$ws.Range("L122").Value = 7066.6671
$ws.Range("M122").Value = -4015.8079
$ws.Range("N122").Value = -11966.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1330.6666
$ws.Range("I126").Value = 746
$ws.Range("K126").Value = 2238
$ws.Range("M126").Value = 232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22867.61
$ws.Range("I132").Value = 27120.975
$ws.Range("J132").Value = 2664.125
$ws.Range("K132").Value = 81362.92499999999
$ws.Range("L132").Value = 7992.375
$ws.Range("M132").Value = -78832.92499999999
$ws.Range("N132").Value = -13052.375
